$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 241, pushing existing rows 241+ down by two.
$ws.Rows.Item(241).Insert()
$ws.Rows.Item(241).Insert()

# --- New row 241 ---
$ws.Cells.Item(241, 1).Value = 7
$ws.Cells.Item(241, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(241, 3).Value = "Ñuble"
$ws.Cells.Item(241, 4).Value = 44466
$ws.Cells.Item(241, 5).Value = 16
$ws.Cells.Item(241, 6).Value = 100112004
$ws.Cells.Item(241, 7).Value = "Cebolla"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "1a (guarda)"
$ws.Cells.Item(241, 10).Value = 300
$ws.Cells.Item(241, 11).Value = 3000
$ws.Cells.Item(241, 12).Value = 3200
$ws.Cells.Item(241, 13).Value = 3100
$ws.Cells.Item(241, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(241, 15).Value = "Región del Maule"
$ws.Cells.Item(241, 16).Value = 194
$ws.Cells.Item(241, 17).Value = 16
$ws.Cells.Item(241, 18).Value = "Hortaliza"

# --- New row 242 ---
$ws.Cells.Item(242, 1).Value = 7
$ws.Cells.Item(242, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(242, 3).Value = "Ñuble"
$ws.Cells.Item(242, 4).Value = 44466
$ws.Cells.Item(242, 5).Value = 16
$ws.Cells.Item(242, 6).Value = 100112004
$ws.Cells.Item(242, 7).Value = "Cebolla"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "2a (guarda)"
$ws.Cells.Item(242, 10).Value = 240
$ws.Cells.Item(242, 11).Value = 2400
$ws.Cells.Item(242, 12).Value = 2500
$ws.Cells.Item(242, 13).Value = 2450
$ws.Cells.Item(242, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 153
$ws.Cells.Item(242, 17).Value = 16
$ws.Cells.Item(242, 18).Value = "Hortaliza"
